$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at the left (A) for the company id, shifting the
# existing "name"/"description" columns one column to the right.
$ws.Columns.Item(1).Insert()

# New header + per-row company id values.
$ws.Range("A1").Value = "companyId"
$companyId = "ef970b3d-5a2b-4b25-9b2e-c3b73d30a5c3"
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = $companyId
}

# Update the (now shifted-to-column-B) "name" header/values.
$ws.Range("B1").Value = "dd"
$ws.Range("B2").Value = "ddss"
$ws.Range("B3").Value = "Human Resourcesdd"
$ws.Range("B4").Value = "markl"
$ws.Range("B5").Value = "Salesss"
$ws.Range("B6").Value = "Customer Supports"
$ws.Range("B7").Value = "IT"

# B2 picks up a small, plain (non-hyperlink) 8pt font.
$ws.Range("B2").ClearFormats()
$ws.Range("B2").Font.Size = 8

# Resize the new column to fit its contents, like the rest of the table.
$ws.Columns.Item(1).AutoFit()

# Leave the selection where the edit finished.
$ws.Range("B2").Select()
